$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "About" sheet: update the version banner + recommended citation
# ---------------------------------------------------------------------
$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: $newVersion"
$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Polosukhinskaya Coal Mine, Russia, M1518, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# ---------------------------------------------------------------------
# 2. "Boundaries and methane sources" sheet: refresh build_version and
#    append the three newly identified point features (rows 15-17)
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# update build_version (column S) for every pre-existing data row (2-14)
for ($r = 2; $r -le 14; $r++) {
    $wsData.Cells.Item($r, 19).Value = $newVersion
}

# Append new rows by cloning the formatting of row 14, then overwrite content
$wsData.Range("A14:T14").Copy()
$wsData.Range("A15:T17").PasteSpecial()
$wsData.Range("E15:E17").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$newRows = @(
    @{
        Row = 15
        B = "M1518.P15"
        C = "ventilation system"
        D = "shaft"
        E = 45751
        G = "Possible ventilation shaft "
        T = "POINT (87.326559 53.920985)"
    },
    @{
        Row = 16
        B = "M1518.P16"
        C = "degasification system"
        D = "drainage station"
        E = 45457
        G = "Possible gas drainage station"
        T = "POINT (87.405355 53.954073)"
    },
    @{
        Row = 17
        B = "M1518.P17"
        C = "degasification system"
        D = "drainage station"
        E = 45457
        G = "Possible gas drainage station"
        T = "POINT (87.406507 53.952673)"
    }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $wsData.Cells.Item($r, 1).Value = "Polosukhinskaya Coal Mine, Russia, M1518"
    $wsData.Cells.Item($r, 2).Value = $nr.B
    $wsData.Cells.Item($r, 3).Value = $nr.C
    $wsData.Cells.Item($r, 4).Value = $nr.D
    $wsData.Cells.Item($r, 5).Value = $nr.E
    $wsData.Cells.Item($r, 6).Value = "Identified visually from Google Earth Pro satellite imagery."
    $wsData.Cells.Item($r, 7).Value = $nr.G
    $wsData.Cells.Item($r, 8).Value = "extracted"
    $wsData.Cells.Item($r, 9).Value = "M1518"
    $wsData.Cells.Item($r, 10).Value = "Rusugolholding LLC [100%]"
    # columns K (Owners (Non-ENG)) and N (GEM Wiki Page (Non-ENG)) stay blank,
    # matching the source row
    $wsData.Cells.Item($r, 12).Value = "Stirakia Holdings Ltd"
    $wsData.Cells.Item($r, 13).Value = "https://www.gem.wiki/Polosukhinskaya_coal_mine"
    $wsData.Cells.Item($r, 15).Value = "Met"
    $wsData.Cells.Item($r, 16).Value = "Polosukhinskaya Coal Mine"
    $wsData.Cells.Item($r, 17).Value = "Russia"
    $wsData.Cells.Item($r, 18).Value = "Nov 26, 2025"
    $wsData.Cells.Item($r, 19).Value = $newVersion
    $wsData.Cells.Item($r, 20).Value = $nr.T
}
